# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the old Late / heading / Outstanding columns one place
# to the right, and make that sheet the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at N; existing N/O/P (Late / heading / Outstanding)
# shift right to O/P/Q.
$ws.Columns("N").Insert()

# Match the column width Excel assigns the newly inserted column (copied
# from its left neighbour, column M).
$ws.Columns("N").ColumnWidth = 10.17

# Make "Repayment schedule" the active sheet/tab and set its selection.
$ws.Activate()
$null = $ws.Range("R9").Select()
